# Generate Report for Handback
# - Updates the Status for the 25575e1d-... (handback) file from
#   "Ready for handoff" to "Handback transform failed" everywhere it is
#   shown (Overview sheet, zh-cn sheet, de-de sheet).
# - Records an Error Detail message in the zh-cn and de-de detail sheets
#   explaining the handback/handoff file name mismatch.
# - Widens the "Error Detail" column so the new, longer message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Update every cell that currently shows the old status text for the
# 25575e1d-... row so the shared status string is replaced consistently.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch in the Error Detail column.
$wsZh.Range("P3").Value = "Handback file name: yko4vndy.u5f is different with handoff file name: 25575e1d-9630-420a-adc5-eb9f14f4bb55.0c6b34ac983fe3e2394d51698cd3afff83fbea3a.zh-cn."
$wsDe.Range("P3").Value = "Handback file name: yko4vndy.u5f is different with handoff file name: 25575e1d-9630-420a-adc5-eb9f14f4bb55.0c6b34ac983fe3e2394d51698cd3afff83fbea3a.de-de."

# Widen the Error Detail column (column P, the 16th column) so the longer
# message text is visible.
$wsZh.Columns.Item(16).ColumnWidth = 40
$wsDe.Columns.Item(16).ColumnWidth = 40
